$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"). Copy H1's formatting
# (bold, bordered, centered header style) onto them first, then set text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-30: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 30; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
